$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 340
$ws1.Range("F5").Value = 3386
$ws1.Range("F6").Value = 2143
$ws1.Range("F7").Value = 410
$ws1.Range("F8").Value = 159
$ws1.Range("F9").Value = 44
$ws1.Range("F10").Value = 30
$ws1.Range("F11").Value = 1251
$ws1.Range("F13").Value = 1469
$ws1.Range("F14").Value = 110

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 340
$ws4.Range("F5").Value = 3386
$ws4.Range("F6").Value = 2143
$ws4.Range("F7").Value = 410
$ws4.Range("F9").Value = 159
$ws4.Range("F10").Value = 44
$ws4.Range("F11").Value = 30
$ws4.Range("F14").Value = 1251
$ws4.Range("F16").Value = 1469
$ws4.Range("F17").Value = 110
